$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.893.94'
$ws.Range("E2").Value = '  -2.74%  '
$ws.Range("D3").Value = '2.618.40'
$ws.Range("E3").Value = '  -1.58%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '576.45'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -3.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '156.43'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.35%  '
$ws.Range("E7").Value = '  +0.00%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.624'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.64%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.120'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -4.82%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '5.83'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  +0.34%  '
$ws.Range("E11").Value = '  -4.07%  '
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '28.28'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -1.19%  '
$ws.Range("D14").Value = '3.092.41'
$ws.Range("E14").Value = '  -1.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000183'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -6.15%  '
$ws.Range("D16").Value = '63.687.11'
$ws.Range("E16").Value = '  -2.88%  '
$ws.Range("D17").Value = '2.618.25'
$ws.Range("E17").Value = '  -2.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.11'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -3.31%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.67'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +3.27%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.57'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -4.16%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '343.67'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -1.62%  '
$ws.Range("E22").Value = '  -0.16%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '67.46'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  -3.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.77'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  -0.68%  '
$ws.Range("E25").Value = '  -1.77%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '595.01'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  +6.15%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.23'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -2.96%  '
$ws.Range("E28").Value = '  -1.62%  '
$ws.Range("B29").Value = 'Kaspa'
$ws.Range("C29").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.161'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -0.79%  '
$ws.Range("B30").Value = 'Binance-PegBSC-USD'
$ws.Range("C30").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.00'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  +0.02%  '
$ws.Range("B31").Value = 'Aptos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.93'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -1.30%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.07'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -2.86%  '
$ws.Range("E33").Value = '  -2.36%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.62'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.48%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.38'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -1.35%  '
$ws.Range("E36").Value = '  -2.93%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '19.78'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -3.15%  '
$ws.Range("E38").Value = '  -0.01%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '154.63'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -0.66%  '
$ws.Range("E40").Value = '  -2.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '41.50'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -2.65%  '
$ws.Range("E43").Value = '  +6.97%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '156.10'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -2.80%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '3.92'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -3.16%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '23.44'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  +3.88%  '
$ws.Range("E47").Value = '  -1.98%  '
$ws.Range("E48").Value = '  +0.03%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.629'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -1.53%  '
$ws.Range("E50").Value = '  -2.77%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.04'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -3.50%  '
